$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells R1:X1 (copy header style from Q1) ---
$headers = @("N", "LTD", "LTM", "LTS", "LND", "LNM", "LNS")
$headerCols = @(18, 19, 20, 21, 22, 23, 24)
for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws.Cells.Item(1, $headerCols[$i]).Value = $headers[$i]
}
$ws.Range("Q1").Copy() | Out-Null
$ws.Range("R1:X1").PasteSpecial(-4122) | Out-Null

# --- Row data: M (text, shifted up one row) + new R..X columns ---
$rows = @(
    @{ Row = 2; M = "100.0516"; R = "48.9676"; S = 41; T = 20; U = 26.7106571498374; V = 2; W = 9; X = 52.46680573167545 },
    @{ Row = 3; M = "51.8971"; R = "48.9561"; S = 41; T = 19; U = 35.70145709855183; V = 2; W = 8; X = 32.11599954935174 },
    @{ Row = 4; M = "52.3443"; R = "48.9543"; S = 41; T = 19; U = 33.16839982050851; V = 2; W = 8; X = 31.46404049253391 },
    @{ Row = 5; M = "51.9766"; R = "48.9406"; S = 41; T = 19; U = 25.11875552282334; V = 2; W = 8; X = 48.65540848184445 },
    @{ Row = 6; M = "51.9616"; R = "48.9406"; S = 41; T = 19; U = 25.11882324071223; V = 2; W = 8; X = 48.65579474154387 },
    @{ Row = 7; M = "51.9184"; R = "48.9354"; S = 41; T = 19; U = 17.98430469970413; V = 2; W = 8; X = 46.76944868692495 },
    @{ Row = 8; M = "51.9414"; R = "48.9354"; S = 41; T = 19; U = 17.98511881471313; V = 2; W = 8; X = 46.76991122247561 },
    @{ Row = 9; M = "52.4805"; R = "48.9475"; S = 41; T = 19; U = 20.63749414482231; V = 2; W = 8; X = 22.82304878752715 },
    @{ Row = 10; M = "52.4905"; R = "48.9475"; S = 41; T = 19; U = 20.63739655082543; V = 2; W = 8; X = 22.82300705611043 },
    @{ Row = 11; M = "52.4765"; R = "48.9635"; S = 41; T = 19; U = 21.64800740920953; V = 2; W = 7; X = 46.99888501452591 },
    @{ Row = 12; M = "52.5015"; R = "48.9635"; S = 41; T = 19; U = 21.64816757838025; V = 2; W = 7; X = 46.99862479207273 },
    @{ Row = 13; M = "52.4975"; R = "48.9635"; S = 41; T = 19; U = 21.64814003256822; V = 2; W = 7; X = 46.99927038701023 },
    @{ Row = 14; M = "53.3178"; R = "48.9598"; S = 41; T = 19; U = 21.38552922561189; V = 2; W = 7; X = 54.66617734017387 },
    @{ Row = 15; M = "53.3108"; R = "48.9598"; S = 41; T = 19; U = 21.38562682239694; V = 2; W = 7; X = 54.66621906000306 },
    @{ Row = 16; M = "52.9321"; R = "48.9521"; S = 41; T = 19; U = 19.61109632877708; V = 2; W = 8; X = 9.614466331483555 },
    @{ Row = 17; M = "52.5073"; R = "48.9513"; S = 41; T = 19; U = 20.44438711033138; V = 2; W = 8; X = 13.19826761875039 },
    @{ Row = 18; M = "52.5698"; R = "48.9568"; S = 41; T = 19; U = 29.30799650146184; V = 2; W = 8; X = 17.46716381549358 },
    @{ Row = 19; M = "52.5588"; R = "48.9568"; S = 41; T = 19; U = 29.30702023417723; V = 2; W = 8; X = 17.46670349646909 },
    @{ Row = 20; M = "52.5648"; R = "48.9568"; S = 41; T = 19; U = 29.30747578263379; V = 2; W = 8; X = 17.46691257702675 },
    @{ Row = 21; M = "52.5718"; R = "48.9598"; S = 41; T = 19; U = 29.30763758129785; V = 2; W = 8; X = 17.46686742840751 }
)

foreach ($d in $rows) {
    $mCell = $ws.Cells.Item($d.Row, 13)
    $mCell.NumberFormat = "@"
    $mCell.Value = $d.M
    $mCell.Style = "Normal"

    $rCell = $ws.Cells.Item($d.Row, 18)
    $rCell.NumberFormat = "@"
    $rCell.Value = $d.R
    $rCell.Style = "Normal"

    $ws.Cells.Item($d.Row, 19).Value = $d.S
    $ws.Cells.Item($d.Row, 20).Value = $d.T
    $ws.Cells.Item($d.Row, 21).Value = $d.U
    $ws.Cells.Item($d.Row, 22).Value = $d.V
    $ws.Cells.Item($d.Row, 23).Value = $d.W
    $ws.Cells.Item($d.Row, 24).Value = $d.X
}

Write-Host "Done applying SURVEYING layout automation edits."